$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "VET_persistence"
$ws.Cells.Item(1,1).Value = "Event"
$ws.Cells.Item(1,2).Value = "Lead Time 1"
$ws.Cells.Item(1,3).Value = "Lead Time 2"
$ws.Cells.Item(1,4).Value = "Lead Time 3"
$ws.Cells.Item(1,5).Value = "Lead Time 4"
$ws.Cells.Item(1,6).Value = "Lead Time 5"
$ws.Cells.Item(1,7).Value = "Lead Time 6"
$ws.Range("A1:G1").Font.Bold = $true
$ws.Range("A1:G1").Borders.LineStyle = 1
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").VerticalAlignment = -4160
$ws.Cells.Item(2,1).Value = "20230520_2235"
$ws.Cells.Item(2,2).Value = 0.34847078528093761
$ws.Cells.Item(2,3).Value = 0.20657729728399099
$ws.Cells.Item(2,4).Value = 0.14283322642220889
$ws.Cells.Item(2,5).Value = 0.10449376677107799
$ws.Cells.Item(2,6).Value = 0.080780644502555149
$ws.Cells.Item(2,7).Value = 0.067002638535033729
$ws.Cells.Item(3,1).Value = "20190320_0005"
$ws.Cells.Item(3,2).Value = 0.42973772801780902
$ws.Cells.Item(3,3).Value = 0.267317091131477
$ws.Cells.Item(3,4).Value = 0.19420543974774271
$ws.Cells.Item(3,5).Value = 0.14615323284923931
$ws.Cells.Item(3,6).Value = 0.10863080162755009
$ws.Cells.Item(3,7).Value = 0.077750884338902898
$ws.Cells.Item(4,1).Value = "20191222_0900"
$ws.Cells.Item(4,2).Value = 0.073506261782311622
$ws.Cells.Item(4,3).Value = 0.04116610517749255
$ws.Cells.Item(4,4).Value = 0.03026140816005644
$ws.Cells.Item(4,5).Value = 0.02697905299759264
$ws.Cells.Item(4,6).Value = 0.040935923285330751
$ws.Cells.Item(4,7).Value = 0.045460779173834343
$ws.Cells.Item(5,1).Value = "20180505_1745"
$ws.Cells.Item(5,2).Value = 0.40653954632564759
$ws.Cells.Item(5,3).Value = 0.25934343340647331
$ws.Cells.Item(5,4).Value = 0.1864329149978555
$ws.Cells.Item(5,5).Value = 0.13548552225614749
$ws.Cells.Item(5,6).Value = 0.096970336628080703
$ws.Cells.Item(5,7).Value = 0.068258077972422601
$ws.Cells.Item(6,1).Value = "20230513_1455"
$ws.Cells.Item(6,2).Value = 0.097807121511178097
$ws.Cells.Item(6,3).Value = 0.063358229707456887
$ws.Cells.Item(6,4).Value = 0.055416255175492632
$ws.Cells.Item(6,5).Value = 0.051360746133378751
$ws.Cells.Item(6,6).Value = 0.043632130485364892
$ws.Cells.Item(6,7).Value = 0.019320874579757969
$ws.Cells.Item(7,1).Value = "20200911_1315"
$ws.Cells.Item(7,2).Value = 0.47796097952526001
$ws.Cells.Item(7,3).Value = 0.28524806061663133
$ws.Cells.Item(7,4).Value = 0.18847598728427589
$ws.Cells.Item(7,5).Value = 0.13118208993216729
$ws.Cells.Item(7,6).Value = 0.089678030164244221
$ws.Cells.Item(7,7).Value = 0.060182814811351293
$ws.Cells.Item(8,1).Value = "20191111_0710"
$ws.Cells.Item(8,2).Value = 0.25306525044456413
$ws.Cells.Item(8,3).Value = 0.16122256394848511
$ws.Cells.Item(8,4).Value = 0.12636935219907469
$ws.Cells.Item(8,5).Value = 0.097234758097743401
$ws.Cells.Item(8,6).Value = 0.091992835936770684
$ws.Cells.Item(8,7).Value = 0.070571704114618605
$ws.Cells.Item(9,1).Value = "20230302_0245"
$ws.Cells.Item(9,2).Value = 0.22571130028517769
$ws.Cells.Item(9,3).Value = 0.10397584918959039
$ws.Cells.Item(9,4).Value = 0.1012622835114096
$ws.Cells.Item(9,5).Value = 0.1141565370792836
$ws.Cells.Item(9,6).Value = 0.124133718072986
$ws.Cells.Item(9,7).Value = 0.1021658779432687
$ws.Cells.Item(10,1).Value = "20190412_1220"
$ws.Cells.Item(10,2).Value = 0.39232696788164839
$ws.Cells.Item(10,3).Value = 0.21706468613028129
$ws.Cells.Item(10,4).Value = 0.13353450666139671
$ws.Cells.Item(10,5).Value = 0.083044379565437976
$ws.Cells.Item(10,6).Value = 0.063112406247336911
$ws.Cells.Item(10,7).Value = 0.047675201910522033
$ws.Cells.Item(11,1).Value = "20200120_1440"
$ws.Cells.Item(11,2).Value = 0.06013213258964837
$ws.Cells.Item(11,3).Value = 0.018439261400206991
$ws.Cells.Item(11,4).Value = 0.024585141814057219
$ws.Cells.Item(11,5).Value = 0.030108384011536211
$ws.Cells.Item(11,6).Value = 0.03096864391086683
$ws.Cells.Item(11,7).Value = 0.029460959935557269
$ws.Cells.Item(12,1).Value = "20230129_2215"
$ws.Cells.Item(12,2).Value = 0.37260787625432867
$ws.Cells.Item(12,3).Value = 0.24412611647619581
$ws.Cells.Item(12,4).Value = 0.1701242340861836
$ws.Cells.Item(12,5).Value = 0.1314104578169023
$ws.Cells.Item(12,6).Value = 0.1077418318120619
$ws.Cells.Item(12,7).Value = 0.089785009060016033
$ws.Cells.Item(13,1).Value = "20181014_0515"
$ws.Cells.Item(13,2).Value = 0.45207980570316991
$ws.Cells.Item(13,3).Value = 0.30822964513480289
$ws.Cells.Item(13,4).Value = 0.2381942629473624
$ws.Cells.Item(13,5).Value = 0.19794567147584791
$ws.Cells.Item(13,6).Value = 0.17471299677294741
$ws.Cells.Item(13,7).Value = 0.15808902178187839

$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "VET_extrapolation"
$ws.Cells.Item(1,1).Value = "Event"
$ws.Cells.Item(1,2).Value = "Lead Time 1"
$ws.Cells.Item(1,3).Value = "Lead Time 2"
$ws.Cells.Item(1,4).Value = "Lead Time 3"
$ws.Cells.Item(1,5).Value = "Lead Time 4"
$ws.Cells.Item(1,6).Value = "Lead Time 5"
$ws.Cells.Item(1,7).Value = "Lead Time 6"
$ws.Range("A1:G1").Font.Bold = $true
$ws.Range("A1:G1").Borders.LineStyle = 1
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").VerticalAlignment = -4160
$ws.Cells.Item(2,1).Value = "20230520_2235"
$ws.Cells.Item(2,2).Value = 0.54305775985406668
$ws.Cells.Item(2,3).Value = 0.37957508616374702
$ws.Cells.Item(2,4).Value = 0.28796804287545841
$ws.Cells.Item(2,5).Value = 0.22922923183821259
$ws.Cells.Item(2,6).Value = 0.1869500041466215
$ws.Cells.Item(2,7).Value = 0.15380936419628571
$ws.Cells.Item(3,1).Value = "20190320_0005"
$ws.Cells.Item(3,2).Value = 0.64834896438464429
$ws.Cells.Item(3,3).Value = 0.49934293669837881
$ws.Cells.Item(3,4).Value = 0.40661988016707601
$ws.Cells.Item(3,5).Value = 0.33760537887481329
$ws.Cells.Item(3,6).Value = 0.27550397585578412
$ws.Cells.Item(3,7).Value = 0.2210150843117773
$ws.Cells.Item(4,1).Value = "20191222_0900"
$ws.Cells.Item(4,2).Value = 0.054545188263365067
$ws.Cells.Item(4,3).Value = 0.018054905652850121
$ws.Cells.Item(4,4).Value = 0.012757172262503849
$ws.Cells.Item(4,5).Value = 0.0087622217476931979
$ws.Cells.Item(4,6).Value = 0.01266875486398447
$ws.Cells.Item(4,7).Value = 0.0050548433379843669
$ws.Cells.Item(5,1).Value = "20180505_1745"
$ws.Cells.Item(5,2).Value = 0.55399305507806029
$ws.Cells.Item(5,3).Value = 0.37286585385665738
$ws.Cells.Item(5,4).Value = 0.26859426514967338
$ws.Cells.Item(5,5).Value = 0.20508435560008889
$ws.Cells.Item(5,6).Value = 0.16316414996814341
$ws.Cells.Item(5,7).Value = 0.13349414280082211
$ws.Cells.Item(6,1).Value = "20230513_1455"
$ws.Cells.Item(6,2).Value = 0.17955047596842921
$ws.Cells.Item(6,3).Value = 0.081450247481258869
$ws.Cells.Item(6,4).Value = 0.037093858319336183
$ws.Cells.Item(6,5).Value = 0.020452783504902209
$ws.Cells.Item(6,6).Value = 0.0074212324812973031
$ws.Cells.Item(6,7).Value = 0.00095825532030025646
$ws.Cells.Item(7,1).Value = "20200911_1315"
$ws.Cells.Item(7,2).Value = 0.55423143957345822
$ws.Cells.Item(7,3).Value = 0.35495217450879302
$ws.Cells.Item(7,4).Value = 0.24771327833511661
$ws.Cells.Item(7,5).Value = 0.1845702414959117
$ws.Cells.Item(7,6).Value = 0.1412584282779879
$ws.Cells.Item(7,7).Value = 0.10982556760443341
$ws.Cells.Item(8,1).Value = "20191111_0710"
$ws.Cells.Item(8,2).Value = 0.49009946341615801
$ws.Cells.Item(8,3).Value = 0.32214717092039402
$ws.Cells.Item(8,4).Value = 0.22669509264979809
$ws.Cells.Item(8,5).Value = 0.16158330651940311
$ws.Cells.Item(8,6).Value = 0.1197440860151482
$ws.Cells.Item(8,7).Value = 0.088538936683520156
$ws.Cells.Item(9,1).Value = "20230302_0245"
$ws.Cells.Item(9,2).Value = 0.40460576621925071
$ws.Cells.Item(9,3).Value = 0.2241751017540399
$ws.Cells.Item(9,4).Value = 0.1213529991790033
$ws.Cells.Item(9,5).Value = 0.068198268881454277
$ws.Cells.Item(9,6).Value = 0.038974797141610673
$ws.Cells.Item(9,7).Value = 0.018948956731904981
$ws.Cells.Item(10,1).Value = "20190412_1220"
$ws.Cells.Item(10,2).Value = 0.57796015856992167
$ws.Cells.Item(10,3).Value = 0.40386024372376028
$ws.Cells.Item(10,4).Value = 0.3057265109464351
$ws.Cells.Item(10,5).Value = 0.23740463622011421
$ws.Cells.Item(10,6).Value = 0.19213212247778769
$ws.Cells.Item(10,7).Value = 0.1647031798347865
$ws.Cells.Item(11,1).Value = "20200120_1440"
$ws.Cells.Item(11,2).Value = 0.1408097749756009
$ws.Cells.Item(11,3).Value = 0.07506164687743655
$ws.Cells.Item(11,4).Value = 0.049848421489331322
$ws.Cells.Item(11,5).Value = 0.02827388925364344
$ws.Cells.Item(11,6).Value = 0.0192763295113141
$ws.Cells.Item(11,7).Value = 0.01180792149726247
$ws.Cells.Item(12,1).Value = "20230129_2215"
$ws.Cells.Item(12,2).Value = 0.46438015810561939
$ws.Cells.Item(12,3).Value = 0.29947066835284553
$ws.Cells.Item(12,4).Value = 0.2059900308073275
$ws.Cells.Item(12,5).Value = 0.13997649954224661
$ws.Cells.Item(12,6).Value = 0.098680005825233869
$ws.Cells.Item(12,7).Value = 0.072896539286550405
$ws.Cells.Item(13,1).Value = "20181014_0515"
$ws.Cells.Item(13,2).Value = 0.54454722325542637
$ws.Cells.Item(13,3).Value = 0.3974189794287305
$ws.Cells.Item(13,4).Value = 0.31183723450903972
$ws.Cells.Item(13,5).Value = 0.25589168231084242
$ws.Cells.Item(13,6).Value = 0.21597086338479049
$ws.Cells.Item(13,7).Value = 0.18507869776337779
